$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "ZAF-2013-E"
$wb.Worksheets.Item(2).Name = "ZAF-2013-X"
